$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D: new value, column E: cell removed entirely (capital structure database update
# dropped historical_growth_net_income_last_5_years for these rows)
$ws.Range("D2").Value = -0.0362
$ws.Range("D3").Value = -0.0362
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()

$ws.Range("G2").Value = -0.3274327122153209
$ws.Range("G3").Value = -0.3274327122153209
$ws.Range("H2").Value = -0.3374741200828157
$ws.Range("H3").Value = -0.3374741200828157
$ws.Range("I2").Value = -0.5403726708074534
$ws.Range("I3").Value = -0.5403726708074534
$ws.Range("J2").Value = -0.5403726708074534
$ws.Range("J3").Value = -0.5403726708074534
$ws.Range("K2").Value = -66.59999999999999
$ws.Range("K3").Value = -66.59999999999999
$ws.Range("L2").Value = -1.37888198757764
$ws.Range("L3").Value = -1.37888198757764

# M,N,P,Q,S unchanged. O/R flip sign of zero between the two rows.
$ws.Range("O2").Value = -0
$ws.Range("R2").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("R3").Value = 0

$ws.Range("U2").Value = 2.81
$ws.Range("U3").Value = 2.81
$ws.Range("V2").Value = 0.02741463414634146
$ws.Range("V3").Value = 0.02741463414634146
$ws.Range("W2").Value = -0.5485996705107083
$ws.Range("W3").Value = -0.5485996705107083
$ws.Range("X2").Value = 0.08339416810924212
$ws.Range("X3").Value = 0.08339416810924212
$ws.Range("Y2").Value = -0.6319938386199505
$ws.Range("Y3").Value = -0.6319938386199505
$ws.Range("Z2").Value = 0.3222793087342363
$ws.Range("Z3").Value = 0.3222793087342363
$ws.Range("AA2").Value = -0.1741509308066991
$ws.Range("AA3").Value = -0.1741509308066991
$ws.Range("AB2").Value = 0.06933017970230612
$ws.Range("AB3").Value = 0.06933017970230612
$ws.Range("AC2").Value = -0.2434811105090053
$ws.Range("AC3").Value = -0.2434811105090053
$ws.Range("AD2").Value = 34.6
$ws.Range("AD3").Value = 34.6
$ws.Range("AF2").Value = 34.6
$ws.Range("AF3").Value = 34.6
$ws.Range("AG2").Value = 31.79
$ws.Range("AG3").Value = 31.79
$ws.Range("AH2").Value = 0.2523705324580598
$ws.Range("AH3").Value = 0.2523705324580598
$ws.Range("AI2").Value = 0.3537832310838446
$ws.Range("AI3").Value = 0.3537832310838446
$ws.Range("AJ2").Value = 0.2367264874525281
$ws.Range("AJ3").Value = 0.2367264874525281
$ws.Range("AK2").Value = 0.3346668070323192
$ws.Range("AK3").Value = 0.3346668070323192
$ws.Range("AL2").Value = 6.2
$ws.Range("AL3").Value = 6.2
$ws.Range("AM2").Value = 6.2
$ws.Range("AM3").Value = 6.2
$ws.Range("AN2").Value = -1.544642857142857
$ws.Range("AN3").Value = -1.544642857142857
$ws.Range("AO2").Value = -4.209677419354839
$ws.Range("AO3").Value = -4.209677419354839
$ws.Range("AP2").Value = -1.419196428571429
$ws.Range("AP3").Value = -1.419196428571429
$ws.Range("AQ2").Value = -4.209677419354839
$ws.Range("AQ3").Value = -4.209677419354839
